# Update the "repaymentstrategy" value on the ProductLoanInput sheet
# from "Mifos style" to "Penalties, Fees, Interest, Principal order",
# reusing the left/top aligned formatting already used by the
# product-name cell (B1) so the lengthy text displays the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B17").Select()
